# CBS_EmployeeBillingUtilization.xlsx - "updated employee billing utilization"
#
# Inserts two new header columns ("company / Initiative utilization " and
# "Utilization on self ") before the existing "EmpID" column, which moves
# from F1 to H1. The two new headers get a bold font + vertically centered
# alignment style; "EmpID" keeps its original bold style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "EmpID" moves from column F to column H, keeping its existing bold style.
$ws.Range("H1").Value = "EmpID"
$ws.Range("H1").Font.Bold = $true

# New headers land in F1/G1 with bold font + vertical-center alignment.
$ws.Range("F1").Value = "company / Initiative utilization "
$ws.Range("G1").Value = "Utilization on self "
$ws.Range("F1:G1").Font.Bold = $true
$ws.Range("F1:G1").VerticalAlignment = -4108

# The header row (and selection) now spans through column H.
[void]$ws.Range("A1:H1").Select()

# Widen/add the affected columns to match the new layout.
$ws.Columns.Item(6).ColumnWidth = 28.333333333333336
$ws.Columns.Item(7).ColumnWidth = 18
$ws.Columns.Item(8).ColumnWidth = 8.833333333333332
